$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7751798033714294
$ws.Range("B1").Value = 2.282429933547974
$ws.Range("D1").Value = 0.8840907216072083
$ws.Range("E1").Value = 0.7185125350952148
